# Auto-generated edit script applying the cryptos.xlsx diff
# (GitHub Actions data refresh: prices + 1h volume % changes, plus a
# Kaspa / WrappedeETH row swap at rows 28-29)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a "plain" numeric-looking string (e.g. "1.00", "0.530")
# must be pre-formatted as Text so Excel keeps them as literal strings instead of
# silently converting them to numbers (which would drop meaningful trailing zeros).
$textForceCells = @(
    "D4", "D5", "D6", "D7", "D8", "D10", "D11", "D12", "D13", "D15", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D35", "D38", "D40", "D42", "D43", "D45", "D47", "D48", "D50"
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "54.501.83"
$ws.Range("E2").Value = "  -3.35%  "
$ws.Range("D3").Value = "2.289.67"
$ws.Range("E3").Value = "  -3.54%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "494.50"
$ws.Range("E5").Value = "  -2.36%  "
$ws.Range("D6").Value = "127.25"
$ws.Range("E6").Value = "  -3.89%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("D8").Value = "0.530"
$ws.Range("E8").Value = "  -2.48%  "
$ws.Range("D9").Value = "2.290.63"
$ws.Range("E9").Value = "  -4.28%  "
$ws.Range("D10").Value = "0.0949"
$ws.Range("E10").Value = "  -2.62%  "
$ws.Range("D11").Value = "0.151"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "0.324"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").Value = "4.63"
$ws.Range("E13").Value = "  -3.77%  "
$ws.Range("D14").Value = "2.698.48"
$ws.Range("E14").Value = "  -3.61%  "
$ws.Range("D15").Value = "21.67"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "54.421.54"
$ws.Range("E16").Value = "  -3.40%  "
$ws.Range("D17").Value = "0.0000130"
$ws.Range("E17").Value = "  -2.99%  "
$ws.Range("D18").Value = "2.300.32"
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("D19").Value = "10.02"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").Value = "4.07"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "303.81"
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("D22").Value = "6.49"
$ws.Range("E22").Value = "  +3.22%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("D24").Value = "5.37"
$ws.Range("E24").Value = "  -3.03%  "
$ws.Range("D25").Value = "63.57"
$ws.Range("E25").Value = "  -3.26%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.93%  "
$ws.Range("D27").Value = "0.374"
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.151"
$ws.Range("E28").Value = "  +2.92%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.394.59"
$ws.Range("E29").Value = "  -3.25%  "
$ws.Range("D30").Value = "7.08"
$ws.Range("E30").Value = "  -2.12%  "
$ws.Range("D31").Value = "169.87"
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("D32").Value = "1.60"
$ws.Range("E32").Value = "  -3.15%  "
$ws.Range("D33").Value = "0.0₃0686"
$ws.Range("E33").Value = "  -4.81%  "
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("E37").Value = "  -2.24%  "
$ws.Range("D38").Value = "17.64"
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").Value = "0.867"
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("D42").Value = "35.52"
$ws.Range("E42").Value = "  -2.48%  "
$ws.Range("D43").Value = "0.375"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("E44").Value = "  -2.01%  "
$ws.Range("D45").Value = "129.92"
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("D47").Value = "4.80"
$ws.Range("E47").Value = "  -3.78%  "
$ws.Range("D48").Value = "0.0894"
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("E49").Value = "  -2.67%  "
$ws.Range("D50").Value = "238.99"
$ws.Range("E50").Value = "  -2.85%  "
$ws.Range("E51").Value = "  -1.01%  "
